$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 104, pushing the existing row 104 (and below) down to row 105.
$ws.Rows.Item(104).Insert()

# Populate the newly inserted row 104 with the new data.
$ws.Range("A104").Value = 4
$ws.Range("B104").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C104").Value = "Los Lagos"
$ws.Range("D104").Value = 44628
$ws.Range("D104").NumberFormat = $ws.Range("D105").NumberFormat
$ws.Range("E104").Value = 10
$ws.Range("F104").Value = 100112052
$ws.Range("G104").Value = "Albahaca"
$ws.Range("H104").Value = "Sin especificar"
$ws.Range("I104").Value = "Primera"
$ws.Range("J104").Value = 180
$ws.Range("K104").Value = 5000
$ws.Range("L104").Value = 6000
$ws.Range("M104").Value = 5500
$ws.Range("N104").Value = '$/docena de matas'
$ws.Range("O104").Value = "Región Metropolitana"
$ws.Range("P104").Value = 917
$ws.Range("Q104").Value = 6
$ws.Range("R104").Value = "Hortaliza"
